$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: K15 changes from the number 112 to the text "150+"
$ws.Range("K15").Value = "150+"

# New row 16: duplicate of the "zeroshot huang combined with own" experiment,
# but with min_sim raised to 0.95 and no quality comment yet
$ws.Range("A16").Value = "zeroshot huang combined with own"
$ws.Range("D16").Value = 3000
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 0.95
$ws.Range("H16").Value = 4000
$ws.Range("I16").Value = "null"
$ws.Range("J16").Value = "yes"
$ws.Range("K16").Value = 112

# Match the yellow highlight formatting used throughout row 15 for the new row
$ws.Range("A16:L16").Interior.Color = 65535

# Update the active selection to I12
$ws.Range("I12").Select()
